$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell -> new text value, derived from the authoritative diff.
# NumberFormat is forced to Text ("@") first so numeric-looking
# strings like "577.04" or "63.131.41" stay text, matching the
# workbook's original inlineStr storage instead of being coerced
# into Excel numbers/dates.
$updates = @(
    @{ Cell = "D2"; Value = '63.131.41' }
    @{ Cell = "E2"; Value = '  +0.23%  ' }
    @{ Cell = "D3"; Value = '2.475.57' }
    @{ Cell = "E3"; Value = '  +0.82%  ' }
    @{ Cell = "E4"; Value = '  +0.02%  ' }
    @{ Cell = "D5"; Value = '577.04' }
    @{ Cell = "E5"; Value = '  +0.58%  ' }
    @{ Cell = "D6"; Value = '146.59' }
    @{ Cell = "E6"; Value = '  +0.38%  ' }
    @{ Cell = "E7"; Value = '  +0.12%  ' }
    @{ Cell = "E8"; Value = '  -0.31%  ' }
    @{ Cell = "D9"; Value = '2.475.26' }
    @{ Cell = "E9"; Value = '  +0.87%  ' }
    @{ Cell = "E10"; Value = '  +0.49%  ' }
    @{ Cell = "E11"; Value = '  +1.68%  ' }
    @{ Cell = "E12"; Value = '  +0.77%  ' }
    @{ Cell = "D13"; Value = '0.353' }
    @{ Cell = "E13"; Value = '  +0.22%  ' }
    @{ Cell = "D14"; Value = '28.63' }
    @{ Cell = "E14"; Value = '  +4.80%  ' }
    @{ Cell = "D15"; Value = '0.0000179' }
    @{ Cell = "E15"; Value = '  +0.95%  ' }
    @{ Cell = "D16"; Value = '2.926.30' }
    @{ Cell = "E16"; Value = '  +2.02%  ' }
    @{ Cell = "D17"; Value = '63.093.97' }
    @{ Cell = "E17"; Value = '  +0.59%  ' }
    @{ Cell = "D18"; Value = '2.475.79' }
    @{ Cell = "E18"; Value = '  +1.21%  ' }
    @{ Cell = "D19"; Value = '8.26' }
    @{ Cell = "E19"; Value = '  +4.22%  ' }
    @{ Cell = "D20"; Value = '11.07' }
    @{ Cell = "E20"; Value = '  +0.90%  ' }
    @{ Cell = "D21"; Value = '329.45' }
    @{ Cell = "E21"; Value = '  +0.23%  ' }
    @{ Cell = "E22"; Value = '  +10.57%  ' }
    @{ Cell = "D23"; Value = '4.13' }
    @{ Cell = "E23"; Value = '  +0.02%  ' }
    @{ Cell = "E24"; Value = '  +0.11%  ' }
    @{ Cell = "D25"; Value = '66.23' }
    @{ Cell = "E25"; Value = '  +0.80%  ' }
    @{ Cell = "D26"; Value = '671.64' }
    @{ Cell = "E26"; Value = '  +5.75%  ' }
    @{ Cell = "D27"; Value = '9.71' }
    @{ Cell = "E27"; Value = '  +14.11%  ' }
    @{ Cell = "D28"; Value = '0.0₃0996' }
    @{ Cell = "E28"; Value = '  +0.66%  ' }
    @{ Cell = "E30"; Value = '  +382.08%  ' }
    @{ Cell = "E31"; Value = '  +2.70%  ' }
    @{ Cell = "D32"; Value = '8.07' }
    @{ Cell = "E32"; Value = '  -1.82%  ' }
    @{ Cell = "E33"; Value = '  +1.06%  ' }
    @{ Cell = "E34"; Value = '  -3.70%  ' }
    @{ Cell = "E35"; Value = '  +3.46%  ' }
    @{ Cell = "E36"; Value = '  +0.00%  ' }
    @{ Cell = "D37"; Value = '4.79' }
    @{ Cell = "E37"; Value = '  +0.71%  ' }
    @{ Cell = "D38"; Value = '5.48' }
    @{ Cell = "E38"; Value = '  +1.30%  ' }
    @{ Cell = "E39"; Value = '  -0.78%  ' }
    @{ Cell = "D40"; Value = '18.80' }
    @{ Cell = "E40"; Value = '  +0.64%  ' }
    @{ Cell = "D41"; Value = '151.79' }
    @{ Cell = "E41"; Value = '  -0.83%  ' }
    @{ Cell = "D42"; Value = '2.74' }
    @{ Cell = "E42"; Value = '  -1.70%  ' }
    @{ Cell = "E43"; Value = '  -0.53%  ' }
    @{ Cell = "E44"; Value = '  +0.00%  ' }
    @{ Cell = "D45"; Value = '0.0₆0314' }
    @{ Cell = "E45"; Value = '  +8.69%  ' }
    @{ Cell = "D46"; Value = '154.65' }
    @{ Cell = "E46"; Value = '  +6.62%  ' }
    @{ Cell = "E47"; Value = '  +15.96%  ' }
    @{ Cell = "E48"; Value = '  +0.11%  ' }
    @{ Cell = "D49"; Value = '20.62' }
    @{ Cell = "E49"; Value = '  +0.53%  ' }
    @{ Cell = "D50"; Value = '0.606' }
    @{ Cell = "E50"; Value = '  +0.67%  ' }
    @{ Cell = "D51"; Value = '0.0513' }
    @{ Cell = "E51"; Value = '  -0.84%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
